$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A (MRN) and column B (Reference Date) values for rows 2-7
$ws.Range("A2").Value = 7506405
$ws.Range("B2").Value = 42628

$ws.Range("A3").Value = 9999997
$ws.Range("B3").Value = 39783

$ws.Range("A4").Value = 3925535
$ws.Range("B4").Value = 40471

$ws.Range("A5").Value = 4224188
$ws.Range("B5").Value = 40998

$ws.Range("A6").Value = 1000007
$ws.Range("B6").Value = 42087

$ws.Range("A7").Value = 2109372
$ws.Range("B7").Value = 40598

# Update the active selection to B11
$ws.Range("B11").Select()
